$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 47657
$ws.Range("E2").Value = 1849
$ws.Range("F2").Value = 1849
$ws.Range("G2").Value = 1172
$ws.Range("H2").Value = 883
$ws.Range("I2").Value = 932
$ws.Range("J2").Value = -48
$ws.Range("K2").Value = 44373
$ws.Range("L2").Value = 27957
$ws.Range("M2").Value = 16416
$ws.Range("N2").Value = 15594
$ws.Range("O2").Value = 822
$ws.Range("P2").Value = 1675
$ws.Range("Q2").Value = 3200
$ws.Range("R2").Value = -2313
$ws.Range("S2").Value = 415
$ws.Range("T2").Value = 2543
$ws.Range("U2").Value = 657
$ws.Range("V2").Value = 21497
$ws.Range("W2").Value = 3.88
$ws.Range("X2").Value = 1.85
$ws.Range("Y2").Value = 6.14
$ws.Range("Z2").Value = 2.04
$ws.Range("AA2").Value = 170.31
$ws.Range("AB2").Value = 816.79
$ws.Range("AC2").Value = 2782
$ws.Range("AD2").Value = 28.97
$ws.Range("AE2").Value = 55894
$ws.Range("AF2").Value = 1.44
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 1.86
$ws.Range("AI2").Value = 45.08
$ws.Range("AJ2").Value = 30467691

# Row 3
$ws.Range("D3").Value = 39345
$ws.Range("E3").Value = 1640
$ws.Range("F3").Value = 1640
$ws.Range("G3").Value = 1696
$ws.Range("H3").Value = 1219
$ws.Range("I3").Value = 1189
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 42106
$ws.Range("L3").Value = 25468
$ws.Range("M3").Value = 16637
$ws.Range("N3").Value = 15565
$ws.Range("O3").Value = 1073
$ws.Range("P3").Value = 1675
$ws.Range("Q3").Value = 4832
$ws.Range("R3").Value = -3913
$ws.Range("S3").Value = -2584
$ws.Range("T3").Value = 3898
$ws.Range("U3").Value = 934
$ws.Range("V3").Value = 19146
$ws.Range("W3").Value = 4.17
$ws.Range("X3").Value = 3.1
$ws.Range("Y3").Value = 7.63
$ws.Range("Z3").Value = 2.82
$ws.Range("AA3").Value = 153.08
$ws.Range("AB3").Value = 844.13
$ws.Range("AC3").Value = 3550
$ws.Range("AD3").Value = 14.67
$ws.Range("AE3").Value = 55791
$ws.Range("AF3").Value = 0.93
$ws.Range("AG3").Value = 800
$ws.Range("AH3").Value = 1.54
$ws.Range("AI3").Value = 18.9
$ws.Range("AJ3").Value = 30467691

# Row 4
$ws.Range("D4").Value = 39704
$ws.Range("E4").Value = 1571
$ws.Range("F4").Value = 1571
$ws.Range("G4").Value = 1129
$ws.Range("H4").Value = 808
$ws.Range("I4").Value = 732
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 45461
$ws.Range("L4").Value = 28182
$ws.Range("M4").Value = 17279
$ws.Range("N4").Value = 16197
$ws.Range("O4").Value = 1082
$ws.Range("P4").Value = 1675
$ws.Range("Q4").Value = 2968
$ws.Range("R4").Value = -2507
$ws.Range("S4").Value = 1036
$ws.Range("T4").Value = 2727
$ws.Range("U4").Value = 240
$ws.Range("V4").Value = 20807
$ws.Range("W4").Value = 3.96
$ws.Range("X4").Value = 2.03
$ws.Range("Y4").Value = 4.61
$ws.Range("Z4").Value = 1.85
$ws.Range("AA4").Value = 163.1
$ws.Range("AB4").Value = 880.5
$ws.Range("AC4").Value = 2187
$ws.Range("AD4").Value = 37.49
$ws.Range("AE4").Value = 58058
$ws.Range("AF4").Value = 1.41
$ws.Range("AG4").Value = 800
$ws.Range("AH4").Value = 0.98
$ws.Range("AI4").Value = 30.17
$ws.Range("AJ4").Value = 30467691

# Row 5
$ws.Range("D5").Value = 50648
$ws.Range("E5").Value = 2626
$ws.Range("F5").Value = 2626
$ws.Range("G5").Value = 2898
$ws.Range("H5").Value = 2176
$ws.Range("I5").Value = 2136
$ws.Range("J5").Value = 41
$ws.Range("K5").Value = 45798
$ws.Range("L5").Value = 26230
$ws.Range("M5").Value = 19568
$ws.Range("N5").Value = 18427
$ws.Range("O5").Value = 1141
$ws.Range("P5").Value = 1675
$ws.Range("Q5").Value = 4210
$ws.Range("R5").Value = -1060
$ws.Range("S5").Value = -2767
$ws.Range("T5").Value = 1060
$ws.Range("U5").Value = 3151
$ws.Range("V5").Value = 18091
$ws.Range("W5").Value = 5.18
$ws.Range("X5").Value = 4.3
$ws.Range("Y5").Value = 12.34
$ws.Range("Z5").Value = 4.77
$ws.Range("AA5").Value = 134.05
$ws.Range("AB5").Value = 996.98
$ws.Range("AC5").Value = 6377
$ws.Range("AD5").Value = 15.6
$ws.Range("AE5").Value = 66048
$ws.Range("AF5").Value = 1.51
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1.01
$ws.Range("AI5").Value = 12.78
$ws.Range("AJ5").Value = 30467691

# Row 6
$ws.Range("D6").Value = 55849
$ws.Range("E6").Value = 5546
$ws.Range("F6").Value = 5546
$ws.Range("G6").Value = 6169
$ws.Range("H6").Value = 5031
$ws.Range("I6").Value = 4912
$ws.Range("K6").Value = 46167
$ws.Range("L6").Value = 22686
$ws.Range("M6").Value = 23481
$ws.Range("N6").Value = 23473
$ws.Range("P6").Value = 1675
$ws.Range("Q6").Value = 4846
$ws.Range("R6").Value = -1907
$ws.Range("S6").Value = -4264
$ws.Range("T6").Value = 1341
$ws.Range("U6").Value = 3505
$ws.Range("V6").Value = 14570
$ws.Range("W6").Value = 9.93
$ws.Range("X6").Value = 9.01
$ws.Range("Y6").Value = 23.45
$ws.Range("Z6").Value = 10.94
$ws.Range("AA6").Value = 96.61
$ws.Range("AB6").Value = 1362.22
$ws.Range("AC6").Value = 14667
$ws.Range("AD6").Value = 5.95
$ws.Range("AE6").Value = 84137
$ws.Range("AF6").Value = 1.04
$ws.Range("AG6").Value = 1350
$ws.Range("AH6").Value = 1.55
$ws.Range("AI6").Value = 7.47
$ws.Range("AJ6").Value = 30467691

# Row 7
$ws.Range("D7").Value = 49809
$ws.Range("E7").Value = 4062
$ws.Range("G7").Value = 4176
$ws.Range("H7").Value = 3291
$ws.Range("I7").Value = 3245
$ws.Range("K7").Value = 47489
$ws.Range("L7").Value = 21103
$ws.Range("M7").Value = 26386
$ws.Range("N7").Value = 26297
$ws.Range("P7").Value = 1672
$ws.Range("Q7").Value = 5147
$ws.Range("R7").Value = -1487
$ws.Range("S7").Value = -2214
$ws.Range("T7").Value = 1946
$ws.Range("U7").Value = 3648
$ws.Range("W7").Value = 8.15
$ws.Range("X7").Value = 6.61
$ws.Range("Y7").Value = 13.04
$ws.Range("Z7").Value = 7.03
$ws.Range("AA7").Value = 79.98
$ws.Range("AC7").Value = 9688
$ws.Range("AD7").Value = 7.86
$ws.Range("AE7").Value = 94259
$ws.Range("AF7").Value = 0.8100000000000001
$ws.Range("AG7").Value = 1654
$ws.Range("AH7").Value = 2.17
$ws.Range("AI7").Value = 15.53

# Row 8
$ws.Range("D8").Value = 50213
$ws.Range("E8").Value = 3266
$ws.Range("G8").Value = 3424
$ws.Range("H8").Value = 2688
$ws.Range("I8").Value = 2599
$ws.Range("K8").Value = 48054
$ws.Range("L8").Value = 19695
$ws.Range("M8").Value = 28359
$ws.Range("N8").Value = 28278
$ws.Range("P8").Value = 1672
$ws.Range("Q8").Value = 3999
$ws.Range("R8").Value = -1475
$ws.Range("S8").Value = -1532
$ws.Range("T8").Value = 1795
$ws.Range("U8").Value = 2567
$ws.Range("W8").Value = 6.5
$ws.Range("X8").Value = 5.35
$ws.Range("Y8").Value = 9.52
$ws.Range("Z8").Value = 5.63
$ws.Range("AA8").Value = 69.45
$ws.Range("AC8").Value = 7759
$ws.Range("AD8").Value = 8.43
$ws.Range("AE8").Value = 101359
$ws.Range("AF8").Value = 0.65
$ws.Range("AG8").Value = 1496
$ws.Range("AH8").Value = 2.29
$ws.Range("AI8").Value = 17.54

# Row 9
$ws.Range("D9").Value = 51299
$ws.Range("E9").Value = 3730
$ws.Range("G9").Value = 3989
$ws.Range("H9").Value = 3085
$ws.Range("I9").Value = 2990
$ws.Range("K9").Value = 49723
$ws.Range("L9").Value = 18923
$ws.Range("M9").Value = 30800
$ws.Range("N9").Value = 30683
$ws.Range("P9").Value = 1672
$ws.Range("Q9").Value = 4152
$ws.Range("R9").Value = -1512
$ws.Range("S9").Value = -1396
$ws.Range("T9").Value = 1749
$ws.Range("U9").Value = 2818
$ws.Range("W9").Value = 7.27
$ws.Range("X9").Value = 6.01
$ws.Range("Y9").Value = 10.14
$ws.Range("Z9").Value = 6.31
$ws.Range("AA9").Value = 61.44
$ws.Range("AC9").Value = 8928
$ws.Range("AD9").Value = 7.33
$ws.Range("AE9").Value = 109979
$ws.Range("AF9").Value = 0.59
$ws.Range("AG9").Value = 1542
$ws.Range("AH9").Value = 2.36
$ws.Range("AI9").Value = 15.72
